# Atualização automática de SAO_LUIZ_GONZAGA.xlsx
#
# - Renames "Paineis DARQ"            -> "PAINEIS DARQ"
# - Renames "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Deletes the "Desarquivamentos Pendentes" sheet entirely
#   (its three exclusive shared strings - "PEDIDOS PENDENTES", "%",
#   "SÃO LUIZ GONZAGA" - are dropped automatically because nothing else
#   in the workbook still references them)

$wb = $excel.ActiveWorkbook

# Suppress the "delete sheet" confirmation prompt so automation doesn't block.
$excel.DisplayAlerts = $false

# Remove the obsolete sheet first (discard the boolean success return value
# so it doesn't leak into the PowerShell output stream).
$null = $wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Rename the remaining sheets to their new (upper-cased / accented) titles.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Keep the original tab ("PAINEIS DARQ") as the active/selected sheet,
# matching the workbook's pre-edit selection state.
$wb.Worksheets.Item("PAINEIS DARQ").Activate()

$excel.DisplayAlerts = $true
